$d = $word.ActiveDocument

# Locate the "Branches " paragraph (it is the only paragraph containing just that text,
# centered and bold, right after the "Sync Fork" section).
$branchesPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq "Branches `r") {
        $branchesPara = $cand
        break
    }
}

if ($branchesPara -eq $null) {
    throw "Could not find the 'Branches' paragraph"
}

# 1) Bold the paragraph mark of the "Branches " paragraph (adds <w:b/> to w:pPr/w:rPr).
$branchesPara.Range.Font.Bold = 1

# 2) Insert the four new explanatory paragraphs (plus one following blank paragraph that
#    matches the paragraph that was already there) right after the "Branches " paragraph.
$insertXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
'<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">If you have the copy of the same code in multiple </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>environment</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> like staging, development then they are called as branches.</w:t></w:r></w:p>' + `
'<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">The main codes are kept in master branch. The same codes are kept in </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>staging ,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> dev environment like that. </w:t></w:r></w:p>' + `
'<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Head represents the latest commit in your repository.</w:t></w:r></w:p>' + `
'<w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>' + `
'<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>' + `
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint = $d.Range($branchesPara.Range.End, $branchesPara.Range.End)
$insertPoint.InsertXML($insertXml)
